# Apply the authoring changes described by the commit "Add files via upload":
#  - Localise sheet names and header labels from English to Catalan.
#  - Fix a data-entry typo on the "Buses" sheet (P value for bus 7).
#  - Update the saved selection (active cell) on both sheets.

$wb = $excel.ActiveWorkbook

$wsBuses    = $wb.Worksheets.Item(1)
$wsTopology = $wb.Worksheets.Item(2)

# --- Rename sheets: Buses -> Busos, Topology -> Topologia ---
$wsBuses.Name    = "Busos"
$wsTopology.Name = "Topologia"

# --- Relabel headers on the "Busos" sheet (Bus / P / Q / V stay the same) ---
$wsBuses.Range("E1").Value = "Tipus"

# --- Relabel headers on the "Topologia" sheet ---
$wsTopology.Range("A1").Value = "Bus inici"
$wsTopology.Range("B1").Value = "Bus fi"

# --- Fix data value: bus 7 P value was -5, should be -0.1 ---
$wsBuses.Range("B9").Value = -0.1

# --- Restore the saved selections on each sheet ---
[void]$wsBuses.Activate()
[void]$wsBuses.Range("D3").Select()

[void]$wsTopology.Activate()
[void]$wsTopology.Range("M12").Select()

# Leave "Busos" as the active/visible tab, matching tabSelected="1" in the target sheet.
[void]$wsBuses.Activate()
